# QAPF with Result Report and Index Showing
# Update the Label/Color/Alpha columns for the "index" data table on the
# "Sheet3" worksheet (the data-node table driving the QAPF plot legend),
# and refresh the sheet selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")
$ws.Activate()

# Resize/reposition the workbook window to match the saved view state.
$win = $excel.ActiveWindow
$win.Left = 3510
$win.Top = 0
$win.Width = 21600
$win.Height = 16200

# Column E (Label) and G (Color) for rows 2-20 are repointed to the new
# "Nodes"/"grey" shared-string values, and column I (Alpha) is dimmed
# from 0.6 to 0.4 for every data row.
for ($row = 2; $row -le 20; $row++) {
    $ws.Cells.Item($row, 5).Value = "Nodes"   # column E - Label
    $ws.Cells.Item($row, 7).Value = "grey"    # column G - Color
    $ws.Cells.Item($row, 9).Value = 0.4       # column I - Alpha
}

# Reflect the new selection recorded for the sheet view.
$ws.Range("E2:E20").Select()
